$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 35720480
$ws.Range("I32").Value = 166668670
$ws.Range("J32").Value = 7335.727
$ws.Range("K32").Value = 166668670
$ws.Range("L32").Value = 7335.727
$ws.Range("M32").Value = -166668344
$ws.Range("N32").Value = -7987.727

$ws.Range("H33").Value = 483.0909
$ws.Range("I33").Value = 483.0909
$ws.Range("K33").Value = 483.0909
$ws.Range("M33").Value = -254.0909

$ws.Range("H92").Value = 45454852
$ws.Range("I92").Value = 52631670
$ws.Range("J92").Value = 1666
$ws.Range("K92").Value = 52631670
$ws.Range("L92").Value = 1666
$ws.Range("M92").Value = -52630422
$ws.Range("N92").Value = -4162

$ws.Range("H98").Value = 572.55554
$ws.Range("I98").Value = 563.8077
$ws.Range("K98").Value = 563.8077
$ws.Range("M98").Value = 934.1923

$ws.Range("H111").Value = 11435.8
$ws.Range("I111").Value = 14497.75
$ws.Range("J111").Value = 10322.363
$ws.Range("K111").Value = 43493.25
$ws.Range("L111").Value = 30967.089
$ws.Range("M111").Value = -40426.25
$ws.Range("N111").Value = -37101.089

$ws.Range("H112").Value = 4102.697
$ws.Range("J112").Value = 4349.6333
$ws.Range("L112").Value = 13048.8999
$ws.Range("N112").Value = -15264.8999

$ws.Range("H118").Value = 512
$ws.Range("I118").Value = 372.77777
$ws.Range("J118").Value = 1138.5
$ws.Range("K118").Value = 1118.33331
$ws.Range("L118").Value = 3415.5
$ws.Range("M118").Value = 538.66669
$ws.Range("N118").Value = -6729.5

$ws.Range("H122").Value = 572.55554
$ws.Range("I122").Value = 563.8077
$ws.Range("K122").Value = 1691.4231
$ws.Range("M122").Value = 758.5769

$ws.Range("H129").Value = 2302
$ws.Range("I129").Value = 633.1667
$ws.Range("K129").Value = 1899.5001
$ws.Range("M129").Value = 3100.4999

$ws.Range("H138").Value = 6116.935
$ws.Range("I138").Value = 2682.4285
$ws.Range("K138").Value = 8047.2855
$ws.Range("M138").Value = -2907.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2844.4167
$ws.Range("I45").Value = 2213.3
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 2213.3
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -1836.3
$ws.Range("N45").Value = -6754

$ws.Range("H74").Value = 1706.2142
$ws.Range("I74").Value = 1288.9
$ws.Range("K74").Value = 1288.9
$ws.Range("M74").Value = -414.9000000000001

$ws.Range("H77").Value = 1706.2142
$ws.Range("I77").Value = 1288.9
$ws.Range("K77").Value = 6444.5
$ws.Range("M77").Value = -2076.5

$ws.Range("H80").Value = 40000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""

$ws.Range("H83").Value = 40000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

$ws.Range("H132").Value = 13587.362
$ws.Range("I132").Value = 14676.488
$ws.Range("K132").Value = 44029.464
$ws.Range("M132").Value = -41499.464

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2026.4546
$ws.Range("I99").Value = 1393.6666
$ws.Range("K99").Value = 1393.6666
$ws.Range("M99").Value = 104.3334

$ws.Range("H107").Value = 2856.7334
$ws.Range("I107").Value = 2346.5
$ws.Range("K107").Value = 2346.5
$ws.Range("M107").Value = -426.5

$ws.Range("H140").Value = 76222.91
$ws.Range("J140").Value = 76295.2
$ws.Range("L140").Value = 76295.2
$ws.Range("N140").Value = -86655.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50002204
$ws.Range("I31").Value = 62501590
$ws.Range("J31").Value = 4648
$ws.Range("K31").Value = 62501590
$ws.Range("L31").Value = 4648
$ws.Range("M31").Value = -62501295
$ws.Range("N31").Value = -5238

$ws.Range("H34").Value = 50002204
$ws.Range("I34").Value = 62501590
$ws.Range("J34").Value = 4648
$ws.Range("K34").Value = 62501590
$ws.Range("L34").Value = 4648
$ws.Range("M34").Value = -62501388
$ws.Range("N34").Value = -5052

$ws.Range("H63").Value = 70000
$ws.Range("J63").Value = 70000
$ws.Range("L63").Value = 70000
$ws.Range("N63").Value = -71372

$ws.Range("H66").Value = 70000
$ws.Range("J66").Value = 70000
$ws.Range("L66").Value = 210000
$ws.Range("N66").Value = -216864

$ws.Range("H99").Value = 28252.2
$ws.Range("I99").Value = 27082.666
$ws.Range("K99").Value = 27082.666
$ws.Range("M99").Value = -25584.666

$ws.Range("H105").Value = 1046.1578
$ws.Range("I105").Value = 826.1429000000001
$ws.Range("K105").Value = 826.1429000000001
$ws.Range("M105").Value = 920.8570999999999

$ws.Range("H126").Value = 28252.2
$ws.Range("I126").Value = 27082.666
$ws.Range("K126").Value = 81247.99800000001
$ws.Range("M126").Value = -78777.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 35715610
$ws.Range("I23").Value = 304.2
$ws.Range("J23").Value = 55557444
$ws.Range("K23").Value = 912.5999999999999
$ws.Range("L23").Value = 166672332
$ws.Range("M23").Value = -677.5999999999999
$ws.Range("N23").Value = -166672802

$ws.Range("H33").Value = 289.66666
$ws.Range("I33").Value = 514.4286
$ws.Range("J33").Value = 93
$ws.Range("K33").Value = 3086.5716
$ws.Range("L33").Value = 558
$ws.Range("M33").Value = -2803.5716
$ws.Range("N33").Value = -1124

$ws.Range("H46").Value = 4999
$ws.Range("J46").Value = 4999
$ws.Range("L46").Value = 14997
$ws.Range("N46").Value = -15179

$ws.Range("H113").Value = 825.0476
$ws.Range("I113").Value = 749.5
$ws.Range("J113").Value = 833
$ws.Range("K113").Value = 2248.5
$ws.Range("L113").Value = 2499
$ws.Range("M113").Value = -78.5
$ws.Range("N113").Value = -6839

$ws.Range("H122").Value = 1074.5264
$ws.Range("J122").Value = 1026.125
$ws.Range("L122").Value = 9235.125
$ws.Range("N122").Value = -14135.125

$ws.Range("H126").Value = 18343.334
$ws.Range("J126").Value = 25000
$ws.Range("L126").Value = 75000
$ws.Range("N126").Value = -84880

$ws.Range("H129").Value = 1370.9048
$ws.Range("J129").Value = 2734.5
$ws.Range("L129").Value = 8203.5
$ws.Range("N129").Value = -18203.5

$ws.Range("H139").Value = 3468.4285
$ws.Range("I139").Value = 2569.75
$ws.Range("J139").Value = 4666.6665
$ws.Range("K139").Value = 7709.25
$ws.Range("L139").Value = 13999.9995
$ws.Range("M139").Value = -2569.25
$ws.Range("N139").Value = -24279.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""

$ws.Range("H102").Value = 8617.809999999999
$ws.Range("I102").Value = 8527.883
$ws.Range("K102").Value = 8527.883
$ws.Range("M102").Value = -6905.883

$ws.Range("H122").Value = 1227145.5
$ws.Range("I122").Value = 1837385
$ws.Range("J122").Value = 6666.6665
$ws.Range("K122").Value = 5512155
$ws.Range("L122").Value = 19999.9995
$ws.Range("M122").Value = -5509705
$ws.Range("N122").Value = -24899.9995

$ws.Range("H132").Value = 2740.2258
$ws.Range("I132").Value = 2239.5518
$ws.Range("K132").Value = 6718.655400000001
$ws.Range("M132").Value = -4188.655400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 702.7619
$ws.Range("J22").Value = 774.8889
$ws.Range("L22").Value = 774.8889
$ws.Range("N22").Value = -1364.8889

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = ""

$ws.Range("H27").Value = 702.7619
$ws.Range("J27").Value = 774.8889
$ws.Range("L27").Value = 774.8889
$ws.Range("N27").Value = -988.8889

$ws.Range("H40").Value = 3292283.8
$ws.Range("I40").Value = 3790100.2
$ws.Range("K40").Value = 3790100.2
$ws.Range("M40").Value = -3789964.2

$ws.Range("H46").Value = 3722.5
$ws.Range("J46").Value = 3815.1738
$ws.Range("L46").Value = 3815.1738
$ws.Range("N46").Value = -4191.1738

$ws.Range("H104").Value = 44164.332
$ws.Range("J104").Value = 44164.332
$ws.Range("L104").Value = 44164.332
$ws.Range("N104").Value = -51152.332

$ws.Range("H132").Value = 3738.4866
$ws.Range("I132").Value = 2680.509
$ws.Range("J132").Value = 6801.0527
$ws.Range("K132").Value = 8041.527
$ws.Range("L132").Value = 20403.1581
$ws.Range("M132").Value = -5511.527
$ws.Range("N132").Value = -25463.1581

$ws.Range("H141").Value = 75316.336
$ws.Range("J141").Value = 75316.336
$ws.Range("L141").Value = 75316.336
$ws.Range("N141").Value = -85676.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 56008.832
$ws.Range("J46").Value = 56008.832
$ws.Range("L46").Value = 56008.832
$ws.Range("N46").Value = -56470.832

$ws.Range("H132").Value = 2022.081
$ws.Range("I132").Value = 1697.9
$ws.Range("J132").Value = 3411.4285
$ws.Range("K132").Value = 5093.700000000001
$ws.Range("L132").Value = 10234.2855
$ws.Range("M132").Value = -2563.700000000001
$ws.Range("N132").Value = -15294.2855

$ws.Range("H134").Value = 56008.832
$ws.Range("J134").Value = 56008.832
$ws.Range("L134").Value = 168026.496
$ws.Range("N134").Value = -173096.496

$ws.Range("H136").Value = 9766.686
$ws.Range("I136").Value = 2760.1875
$ws.Range("J136").Value = 11842.686
$ws.Range("K136").Value = 8280.5625
$ws.Range("L136").Value = 35528.058
$ws.Range("M136").Value = -5730.5625
$ws.Range("N136").Value = -40628.058
